$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump version 1.0.0 -> 1.1.0
$ws.Range("B3").Value = "1.1.0"

# Update date
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
